$wb = $excel.ActiveWorkbook

$wsAstronauta = $wb.Worksheets.Item("Astronauta")
$wsAstronauta.Activate()
$wsAstronauta.Range("E2").Select()

$ws = $wb.Worksheets.Item("Mago")

$ws.Range("E8").Formula = "=""0.5"""
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial(-4163)
$ws.Range("E10").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("E13").Value = 0
$ws.Range("E15").Formula = "=""0.5"""
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial(-4163)
$ws.Range("E16").Value = 0
$ws.Range("E17").Value = 1
$ws.Range("E18").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("E22").Value = 1

$ws.Activate()
$ws.Range("E12").Select()
